$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.778.46"
$ws.Range("E2").Value = "  +9.33%  "

$ws.Range("D3").Value = "2.595.02"
$ws.Range("E3").Value = "  +6.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "186.52"
$ws.Range("E5").Value = "  +16.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "582.50"
$ws.Range("E6").Value = "  +4.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.211"
$ws.Range("E8").Value = "  +29.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  +5.05%  "

$ws.Range("D10").Value = "2.590.44"
$ws.Range("E10").Value = "  +6.73%  "

$ws.Range("E11").Value = "  -0.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.361"
$ws.Range("E12").Value = "  +9.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.82"
$ws.Range("E13").Value = "  +4.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000194"
$ws.Range("E14").Value = "  +12.02%  "

$ws.Range("D15").Value = "74.576.98"
$ws.Range("E15").Value = "  +9.25%  "

$ws.Range("D16").Value = "3.057.89"
$ws.Range("E16").Value = "  +6.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.37"
$ws.Range("E17").Value = "  +14.70%  "

$ws.Range("D18").Value = "2.603.52"
$ws.Range("E18").Value = "  +7.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.49"
$ws.Range("E19").Value = "  +23.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.75"
$ws.Range("E20").Value = "  +12.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.46"
$ws.Range("E21").Value = "  +12.36%  "

$ws.Range("E22").Value = "  +20.52%  "

$ws.Range("E23").Value = "  +6.89%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.84"
$ws.Range("E25").Value = "  +4.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.21"
$ws.Range("E26").Value = "  +14.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("E27").Value = "  +12.59%  "

$ws.Range("D28").Value = "2.728.28"
$ws.Range("E28").Value = "  +7.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").Value = "0.0₃0945"
$ws.Range("E30").Value = "  +16.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  +11.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "504.37"
$ws.Range("E32").Value = "  +18.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  +16.99%  "

$ws.Range("E34").Value = "  +6.57%  "

$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  +14.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.96"
$ws.Range("E37").Value = "  +0.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.24"
$ws.Range("E38").Value = "  +7.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.39"
$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.99"
$ws.Range("E41").Value = "  +15.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.69"
$ws.Range("E42").Value = "  +13.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.321"
$ws.Range("E43").Value = "  +8.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.45"
$ws.Range("E44").Value = "  +22.71%  "

$ws.Range("E45").Value = "  +4.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  +8.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.87"
$ws.Range("E47").Value = "  +12.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0815"
$ws.Range("E48").Value = "  +14.37%  "

$ws.Range("E49").Value = "  +8.29%  "

$ws.Range("E50").Value = "  +8.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.584"
$ws.Range("E51").Value = "  +5.07%  "
